$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.93
$ws.Range("B10").Value = 6.784999999999999
$ws.Range("B12").Value = 6.444
$ws.Range("E13").Value = 12.583
$ws.Range("B18").Value = 6.873
$ws.Range("B25").Value = 6.629
